$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 5351.066666666667, 5697, 4691, 0.1173476060231527),
    @(1, 5241.366666666667, 5584, 4762, 0.1177621046702067),
    @(2, 5180.066666666667, 5603, 4199, 0.1190128167470296),
    @(3, 5670.966666666666, 5963, 5325, 0.1182988484700521),
    @(4, 5042.566666666667, 5358, 4432, 0.1239565928777059),
    @(5, 5010.633333333333, 5413, 4378, 0.1186003843943278),
    @(6, 5529.1,            5984, 4872, 0.1214144468307495),
    @(7, 5383.833333333333, 5789, 4916, 0.121527640024821),
    @(8, 5328.633333333333, 5577, 4668, 0.1228060881296794),
    @(9, 5333.666666666667, 5927, 4877, 0.1205890417098999)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}
